$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrige los valores del modelo para Barcelona (fila 2) y Real Madrid (fila 3)
$ws.Cells.Item(2, 8).Value = 676
$ws.Cells.Item(2, 9).Value = 61
$ws.Cells.Item(2, 10).Value = 266

$ws.Cells.Item(3, 8).Value = 634
$ws.Cells.Item(3, 9).Value = 57
$ws.Cells.Item(3, 10).Value = 243

# Agrega una nueva fila de estadisticas para el Psg
$ws.Cells.Item(4, 1).Value = "Psg"
$ws.Cells.Item(4, 2).Value = 34
$ws.Cells.Item(4, 3).Value = 26
$ws.Cells.Item(4, 4).Value = 6
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 92
$ws.Cells.Item(4, 7).Value = 35
$ws.Cells.Item(4, 8).Value = 636
$ws.Cells.Item(4, 9).Value = 41
$ws.Cells.Item(4, 10).Value = 218

# Actualiza la celda seleccionada
$ws.Range("A5").Select()
